$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 14 de Julio de 2020 a las 22:45"

# --- Estados Unidos (row 4) : updated stats, same rank ---
$ws.Range("B4").Value = 3525355
$ws.Range("C4").Value = 45872
$ws.Range("D4").Value = 1581255
$ws.Range("E4").Value = 1805157
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 696
$ws.Range("H4").Value = 138943

# --- India (row 6) : updated stats, same rank ---
$ws.Range("B6").Value = 937487
$ws.Range("C6").Value = 29842
$ws.Range("D6").Value = 593080
$ws.Range("E6").Value = 320092
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 588
$ws.Range("H6").Value = 24315

# --- China / Egipto swap rank (rows 26-27): Egipto overtakes China ---
$ws.Range("A26").Value = "Egipto"
$ws.Range("B26").Value = 83930
$ws.Range("C26").Value = 929
$ws.Range("D26").Value = 25544
$ws.Range("E26").Value = 54378
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = 73
$ws.Range("H26").Value = 4008

$ws.Range("A27").Value = "China"
$ws.Range("B27").Value = 83605
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 78674
$ws.Range("E27").Value = 297
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 0
$ws.Range("H27").Value = 4634

# --- Emiratos Arabes Unidos (row 38) : updated stats, same rank ---
$ws.Range("B38").Value = 55573
$ws.Range("C38").Value = 375
$ws.Range("D38").Value = 46025
$ws.Range("E38").Value = 9213
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 335

# --- Suiza (row 52) : partial stat update ---
$ws.Range("D52").Value = 29800
$ws.Range("E52").Value = 1248

# --- Uzbekistan (row 67) : updated stats, same rank ---
$ws.Range("B67").Value = 14085
$ws.Range("C67").Value = 494
$ws.Range("D67").Value = 8327
$ws.Range("E67").Value = 5691
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 3
$ws.Range("H67").Value = 67

# --- Costa de Marfil (row 71) : updated stats, same rank ---
$ws.Range("B71").Value = 13037
$ws.Range("C71").Value = 165
$ws.Range("D71").Value = 6908
$ws.Range("E71").Value = 6042
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 3
$ws.Range("H71").Value = 87

# --- Sudan (row 74) : updated stats, same rank ---
$ws.Range("B74").Value = 10417
$ws.Range("C74").Value = 101
$ws.Range("D74").Value = 5579
$ws.Range("E74").Value = 4179
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 659

# --- Guinea (row 92) : updated stats, same rank ---
$ws.Range("B92").Value = 6200
$ws.Range("C92").Value = 59
$ws.Range("D92").Value = 4951
$ws.Range("E92").Value = 1211
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = 38

# --- Malaui moves up ahead of Libano & Cuba (rows 111-113) ---
$ws.Range("A111").Value = "Malaui"
$ws.Range("B111").Value = 2497
$ws.Range("C111").Value = 67
$ws.Range("D111").Value = 795
$ws.Range("E111").Value = 1662
$ws.Range("F111").Value = 0
$ws.Range("G111").Value = 1
$ws.Range("H111").Value = 40

$ws.Range("A112").Value = "Libano"
$ws.Range("B112").Value = 2451
$ws.Range("C112").Value = 32
$ws.Range("D112").Value = 1452
$ws.Range("E112").Value = 962
$ws.Range("F112").Value = 0
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 37

$ws.Range("A113").Value = "Cuba"
$ws.Range("B113").Value = 2432
$ws.Range("C113").Value = 4
$ws.Range("D113").Value = 2275
$ws.Range("E113").Value = 70
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 87

# --- Ruanda overtakes Benin (rows 131-132) ---
$ws.Range("A131").Value = "Ruanda"
$ws.Range("B131").Value = 1416
$ws.Range("C131").Value = 38
$ws.Range("D131").Value = 737
$ws.Range("E131").Value = 675
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 0
$ws.Range("H131").Value = 4

$ws.Range("A132").Value = "Benin"
$ws.Range("B132").Value = 1378
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 557
$ws.Range("E132").Value = 795
$ws.Range("F132").Value = 0
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 26

# --- Togo (row 152) : partial stat update ---
$ws.Range("B152").Value = 731
$ws.Range("C152").Value = 10
$ws.Range("D152").Value = 528
$ws.Range("E152").Value = 188

# --- Angola (row 157) : partial stat update ---
$ws.Range("B157").Value = 541
$ws.Range("C157").Value = 16
$ws.Range("E157").Value = 397
